$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 100
$ws_ALC.Range("H100").Value2 = 1768
$ws_ALC.Range("I100").Value2 = 1711.8182
$ws_ALC.Range("J100").Value2 = 1829.8
$ws_ALC.Range("K100").Value2 = 1711.8182
$ws_ALC.Range("L100").Value2 = 1829.8
$ws_ALC.Range("M100").Value2 = -1170.8182
$ws_ALC.Range("N100").Value2 = -2911.8

# ALC row 116
$ws_ALC.Range("H116").Value2 = 7619.4644
$ws_ALC.Range("I116").Value2 = 5112.222
$ws_ALC.Range("J116").Value2 = 8807.105
$ws_ALC.Range("K116").Value2 = 5112.222
$ws_ALC.Range("L116").Value2 = 8807.105
$ws_ALC.Range("M116").Value2 = -1670.222
$ws_ALC.Range("N116").Value2 = -15691.105

# ALC row 121
$ws_ALC.Range("H121").Value2 = 1588201.2
$ws_ALC.Range("I121").Value2 = 265
$ws_ALC.Range("J121").Value2 = 1852857.2
$ws_ALC.Range("K121").Value2 = 795
$ws_ALC.Range("L121").Value2 = 5558571.6
$ws_ALC.Range("M121").Value2 = 952
$ws_ALC.Range("N121").Value2 = -5562065.6

# ALC row 134
$ws_ALC.Range("H134").Value2 = 53166.668
$ws_ALC.Range("I134").Value2 = 0
$ws_ALC.Range("J134").Value2 = 53166.668
$ws_ALC.Range("K134").Value2 = 0
$ws_ALC.Range("L134").Value2 = 53166.668
$ws_ALC.Range("N134").Value2 = -63306.668

# ARM row 32
$ws_ARM.Range("H32").Value2 = 32766.225
$ws_ARM.Range("I32").Value2 = 31844.38
$ws_ARM.Range("J32").Value2 = 49666.668
$ws_ARM.Range("K32").Value2 = 31844.38
$ws_ARM.Range("L32").Value2 = 49666.668
$ws_ARM.Range("M32").Value2 = -31557.38
$ws_ARM.Range("N32").Value2 = -50240.668

# ARM row 61
$ws_ARM.Range("H61").Value2 = 3124.1304
$ws_ARM.Range("I61").Value2 = 3063.9048
$ws_ARM.Range("J61").Value2 = 3756.5
$ws_ARM.Range("K61").Value2 = 3063.9048
$ws_ARM.Range("L61").Value2 = 3756.5
$ws_ARM.Range("M61").Value2 = -2851.9048
$ws_ARM.Range("N61").Value2 = -4180.5

# ARM row 122
$ws_ARM.Range("H122").Value2 = 3016.4285
$ws_ARM.Range("I122").Value2 = 3808.75
$ws_ARM.Range("J122").Value2 = 1960
$ws_ARM.Range("K122").Value2 = 11426.25
$ws_ARM.Range("L122").Value2 = 5880
$ws_ARM.Range("M122").Value2 = -8976.25
$ws_ARM.Range("N122").Value2 = -10780

# ARM row 136
$ws_ARM.Range("H136").Value2 = 3124.1304
$ws_ARM.Range("I136").Value2 = 3063.9048
$ws_ARM.Range("J136").Value2 = 3756.5
$ws_ARM.Range("K136").Value2 = 9191.714399999999
$ws_ARM.Range("L136").Value2 = 11269.5
$ws_ARM.Range("M136").Value2 = -6641.714399999999
$ws_ARM.Range("N136").Value2 = -16369.5

# BSM row 134
$ws_BSM.Range("H134").Value2 = 2273.818
$ws_BSM.Range("I134").Value2 = 1510.0975
$ws_BSM.Range("J134").Value2 = 4510.4287
$ws_BSM.Range("K134").Value2 = 4530.2925
$ws_BSM.Range("L134").Value2 = 13531.2861
$ws_BSM.Range("M134").Value2 = -1995.2925
$ws_BSM.Range("N134").Value2 = -18601.2861

# CRP row 21
$ws_CRP.Range("H21").Value2 = 10000
$ws_CRP.Range("I21").Value2 = 0
$ws_CRP.Range("J21").Value2 = 10000
$ws_CRP.Range("K21").Value2 = 0
$ws_CRP.Range("L21").Value2 = 10000
$ws_CRP.Range("N21").Value2 = -10470

# CRP row 31
$ws_CRP.Range("H31").Value2 = 6081.6665
$ws_CRP.Range("I31").Value2 = 2823.4546
$ws_CRP.Range("J31").Value2 = 7050.324
$ws_CRP.Range("K31").Value2 = 2823.4546
$ws_CRP.Range("L31").Value2 = 7050.324
$ws_CRP.Range("M31").Value2 = -2528.4546
$ws_CRP.Range("N31").Value2 = -7640.324

# CRP row 34
$ws_CRP.Range("H34").Value2 = 6081.6665
$ws_CRP.Range("I34").Value2 = 2823.4546
$ws_CRP.Range("J34").Value2 = 7050.324
$ws_CRP.Range("K34").Value2 = 2823.4546
$ws_CRP.Range("L34").Value2 = 7050.324
$ws_CRP.Range("M34").Value2 = -2621.4546
$ws_CRP.Range("N34").Value2 = -7454.324

# CRP row 41
$ws_CRP.Range("H41").Value2 = 25583
$ws_CRP.Range("I41").Value2 = 6000
$ws_CRP.Range("J41").Value2 = 35374.5
$ws_CRP.Range("K41").Value2 = 6000
$ws_CRP.Range("L41").Value2 = 35374.5
$ws_CRP.Range("M41").Value2 = -5572
$ws_CRP.Range("N41").Value2 = -36230.5

# CRP row 56
$ws_CRP.Range("H56").Value2 = 15000
$ws_CRP.Range("I56").Value2 = 15000
$ws_CRP.Range("J56").Value2 = 0
$ws_CRP.Range("K56").Value2 = 15000
$ws_CRP.Range("L56").Value2 = 0
$ws_CRP.Range("M56").Value2 = -14155

# CRP row 58
$ws_CRP.Range("H58").Value2 = 1981.1608
$ws_CRP.Range("I58").Value2 = 1700.7872
$ws_CRP.Range("J58").Value2 = 3445.3333
$ws_CRP.Range("K58").Value2 = 1700.7872
$ws_CRP.Range("L58").Value2 = 3445.3333
$ws_CRP.Range("M58").Value2 = -1497.7872
$ws_CRP.Range("N58").Value2 = -3851.3333

# CRP row 63
$ws_CRP.Range("H63").Value2 = 29995
$ws_CRP.Range("I63").Value2 = 0
$ws_CRP.Range("J63").Value2 = 29995
$ws_CRP.Range("K63").Value2 = 0
$ws_CRP.Range("L63").Value2 = 29995
$ws_CRP.Range("N63").Value2 = -31367

# CRP row 66
$ws_CRP.Range("H66").Value2 = 29995
$ws_CRP.Range("I66").Value2 = 0
$ws_CRP.Range("J66").Value2 = 29995
$ws_CRP.Range("K66").Value2 = 0
$ws_CRP.Range("L66").Value2 = 89985
$ws_CRP.Range("N66").Value2 = -96849

# CRP row 99
$ws_CRP.Range("H99").Value2 = 2221
$ws_CRP.Range("I99").Value2 = 2179.111
$ws_CRP.Range("J99").Value2 = 2250
$ws_CRP.Range("K99").Value2 = 2179.111
$ws_CRP.Range("L99").Value2 = 2250
$ws_CRP.Range("M99").Value2 = -681.1109999999999
$ws_CRP.Range("N99").Value2 = -5246

# CRP row 126
$ws_CRP.Range("H126").Value2 = 2221
$ws_CRP.Range("I126").Value2 = 2179.111
$ws_CRP.Range("J126").Value2 = 2250
$ws_CRP.Range("K126").Value2 = 6537.333
$ws_CRP.Range("L126").Value2 = 6750
$ws_CRP.Range("M126").Value2 = -4067.333
$ws_CRP.Range("N126").Value2 = -11690

# CRP row 132
$ws_CRP.Range("H132").Value2 = 38993.973
$ws_CRP.Range("I132").Value2 = 1485.4138
$ws_CRP.Range("J132").Value2 = 159854.89
$ws_CRP.Range("K132").Value2 = 4456.2414
$ws_CRP.Range("L132").Value2 = 479564.67
$ws_CRP.Range("M132").Value2 = -1926.2414
$ws_CRP.Range("N132").Value2 = -484624.67

# CRP row 136
$ws_CRP.Range("H136").Value2 = 1981.1608
$ws_CRP.Range("I136").Value2 = 1700.7872
$ws_CRP.Range("J136").Value2 = 3445.3333
$ws_CRP.Range("K136").Value2 = 5102.3616
$ws_CRP.Range("L136").Value2 = 10335.9999
$ws_CRP.Range("M136").Value2 = -2552.3616
$ws_CRP.Range("N136").Value2 = -15435.9999

# CRP row 141
$ws_CRP.Range("H141").Value2 = 4000
$ws_CRP.Range("I141").Value2 = 4000
$ws_CRP.Range("J141").Value2 = 0
$ws_CRP.Range("K141").Value2 = 4000
$ws_CRP.Range("L141").Value2 = 0
$ws_CRP.Range("M141").ClearContents() | Out-Null
$ws_CRP.Range("N141").Value2 = 1180

# GSM row 68
$ws_GSM.Range("H68").Value2 = 30000
$ws_GSM.Range("I68").Value2 = 0
$ws_GSM.Range("J68").Value2 = 30000
$ws_GSM.Range("K68").Value2 = 0
$ws_GSM.Range("L68").Value2 = 30000
$ws_GSM.Range("N68").Value2 = -31622

# GSM row 71
$ws_GSM.Range("H71").Value2 = 30000
$ws_GSM.Range("I71").Value2 = 0
$ws_GSM.Range("J71").Value2 = 30000
$ws_GSM.Range("K71").Value2 = 0
$ws_GSM.Range("L71").Value2 = 90000
$ws_GSM.Range("N71").Value2 = -98112

# GSM row 109
$ws_GSM.Range("H109").Value2 = 40285
$ws_GSM.Range("I109").Value2 = 0
$ws_GSM.Range("J109").Value2 = 40285
$ws_GSM.Range("K109").Value2 = 0
$ws_GSM.Range("L109").Value2 = 40285
$ws_GSM.Range("N109").Value2 = -42365

# GSM row 122
$ws_GSM.Range("H122").Value2 = 2099.1428
$ws_GSM.Range("I122").Value2 = 2545
$ws_GSM.Range("J122").Value2 = 1920.8
$ws_GSM.Range("K122").Value2 = 7635
$ws_GSM.Range("L122").Value2 = 5762.4
$ws_GSM.Range("M122").Value2 = -5185
$ws_GSM.Range("N122").Value2 = -10662.4

# GSM row 123
$ws_GSM.Range("H123").Value2 = 18660
$ws_GSM.Range("I123").Value2 = 0
$ws_GSM.Range("J123").Value2 = 18660
$ws_GSM.Range("K123").Value2 = 0
$ws_GSM.Range("L123").Value2 = 18660
$ws_GSM.Range("N123").Value2 = -23560

# GSM row 126
$ws_GSM.Range("H126").Value2 = 4616.923
$ws_GSM.Range("I126").Value2 = 7564.4443
$ws_GSM.Range("J126").Value2 = 2090.476
$ws_GSM.Range("K126").Value2 = 22693.3329
$ws_GSM.Range("L126").Value2 = 6271.428
$ws_GSM.Range("M126").Value2 = -20223.3329
$ws_GSM.Range("N126").Value2 = -11211.428

# GSM row 132
$ws_GSM.Range("H132").Value2 = 2037.5471
$ws_GSM.Range("I132").Value2 = 1625.4286
$ws_GSM.Range("J132").Value2 = 3611.0908
$ws_GSM.Range("K132").Value2 = 4876.2858
$ws_GSM.Range("L132").Value2 = 10833.2724
$ws_GSM.Range("M132").Value2 = -2346.2858
$ws_GSM.Range("N132").Value2 = -15893.2724

# LTW row 22
$ws_LTW.Range("H22").Value2 = 1173.5
$ws_LTW.Range("I22").Value2 = 1217.8
$ws_LTW.Range("J22").Value2 = 1099.6666
$ws_LTW.Range("K22").Value2 = 1217.8
$ws_LTW.Range("L22").Value2 = 1099.6666
$ws_LTW.Range("M22").Value2 = -922.8
$ws_LTW.Range("N22").Value2 = -1689.6666

# LTW row 27
$ws_LTW.Range("H27").Value2 = 1173.5
$ws_LTW.Range("I27").Value2 = 1217.8
$ws_LTW.Range("J27").Value2 = 1099.6666
$ws_LTW.Range("K27").Value2 = 1217.8
$ws_LTW.Range("L27").Value2 = 1099.6666
$ws_LTW.Range("M27").Value2 = -1110.8
$ws_LTW.Range("N27").Value2 = -1313.6666

# LTW row 55
$ws_LTW.Range("H55").Value2 = 673.3913
$ws_LTW.Range("I55").Value2 = 717.7273
$ws_LTW.Range("J55").Value2 = 632.75
$ws_LTW.Range("K55").Value2 = 717.7273
$ws_LTW.Range("L55").Value2 = 632.75
$ws_LTW.Range("M55").Value2 = -544.7273
$ws_LTW.Range("N55").Value2 = -978.75

# LTW row 132
$ws_LTW.Range("H132").Value2 = 3246.551
$ws_LTW.Range("I132").Value2 = 3349.923
$ws_LTW.Range("J132").Value2 = 3129.6956
$ws_LTW.Range("K132").Value2 = 10049.769
$ws_LTW.Range("L132").Value2 = 9389.086800000001
$ws_LTW.Range("M132").Value2 = -7519.769
$ws_LTW.Range("N132").Value2 = -14449.0868

# LTW row 140
$ws_LTW.Range("H140").Value2 = 38540.715
$ws_LTW.Range("I140").Value2 = 0
$ws_LTW.Range("J140").Value2 = 38540.715
$ws_LTW.Range("K140").Value2 = 0
$ws_LTW.Range("L140").Value2 = 38540.715
$ws_LTW.Range("N140").Value2 = -48900.715

# WVR row 122
$ws_WVR.Range("H122").Value2 = 29365934
$ws_WVR.Range("I122").Value2 = 37755900
$ws_WVR.Range("J122").Value2 = 1050
$ws_WVR.Range("K122").Value2 = 113267700
$ws_WVR.Range("L122").Value2 = 3150
$ws_WVR.Range("M122").Value2 = -113265250
$ws_WVR.Range("N122").Value2 = -8050

# WVR row 136
$ws_WVR.Range("H136").Value2 = 16956.682
$ws_WVR.Range("I136").Value2 = 39543.117
$ws_WVR.Range("J136").Value2 = 2275.5
$ws_WVR.Range("K136").Value2 = 118629.351
$ws_WVR.Range("L136").Value2 = 6826.5
$ws_WVR.Range("M136").Value2 = -116079.351
$ws_WVR.Range("N136").Value2 = -11926.5
